$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Report 1")
$ws2 = $wb.Worksheets.Item("Report 2")
$ws3 = $wb.Worksheets.Item("Report 3")

$ws1.Range("B12").Value = "source"
$ws1.Range("B13").Value = "target"
$ws1.Range("C12").Value = "input"
$ws1.Range("C13").Value = "output"

$ws1.Activate()
$ws1.Range("G2").Select()

$ws2.Activate()
$ws2.Range("B15").Select()

$ws3.Activate()
$ws3.Range("A11").Select()

$ws1.Activate()
